$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 488.86206
$ws.Range("I33").Value = 576.6818
$ws.Range("J33").Value = 212.85715
$ws.Range("K33").Value = 576.6818
$ws.Range("L33").Value = 212.85715
$ws.Range("M33").Value = -347.6818
$ws.Range("N33").Value = -670.85715
$ws.Range("H41").Value = 191.9375
$ws.Range("I41").Value = 294.4
$ws.Range("J41").Value = 145.36363
$ws.Range("K41").Value = 294.4
$ws.Range("L41").Value = 145.36363
$ws.Range("M41").Value = 145.6
$ws.Range("N41").Value = -1025.36363
$ws.Range("H86").Value = 2577.5
$ws.Range("I86").Value = 2330.5557
$ws.Range("K86").Value = 2330.5557
$ws.Range("M86").Value = -1207.5557
$ws.Range("H88").Value = 2468.9355
$ws.Range("I88").Value = 1345.625
$ws.Range("J88").Value = 2859.652
$ws.Range("K88").Value = 1345.625
$ws.Range("L88").Value = 2859.652
$ws.Range("M88").Value = -939.625
$ws.Range("N88").Value = -3671.652
$ws.Range("H89").Value = 2577.5
$ws.Range("I89").Value = 2330.5557
$ws.Range("K89").Value = 11652.7785
$ws.Range("M89").Value = -6036.7785
$ws.Range("H91").Value = 2468.9355
$ws.Range("I91").Value = 1345.625
$ws.Range("J91").Value = 2859.652
$ws.Range("K91").Value = 1345.625
$ws.Range("L91").Value = 2859.652
$ws.Range("M91").Value = 58.375
$ws.Range("N91").Value = -5667.652
$ws.Range("H137").Value = 1697.2667
$ws.Range("I137").Value = 1231.561
$ws.Range("K137").Value = 3694.683
$ws.Range("M137").Value = -1144.683

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17976.219
$ws.Range("I32").Value = 8364.454
$ws.Range("J32").Value = 57624.75
$ws.Range("K32").Value = 8364.454
$ws.Range("L32").Value = 57624.75
$ws.Range("M32").Value = -8077.454
$ws.Range("N32").Value = -58198.75
$ws.Range("H74").Value = 5244.9585
$ws.Range("I74").Value = 637.0952
$ws.Range("K74").Value = 637.0952
$ws.Range("M74").Value = 236.9048
$ws.Range("H77").Value = 5244.9585
$ws.Range("I77").Value = 637.0952
$ws.Range("K77").Value = 3185.476
$ws.Range("M77").Value = 1182.524
$ws.Range("H88").Value = 3711
$ws.Range("J88").Value = 4663.364
$ws.Range("L88").Value = 4663.364
$ws.Range("N88").Value = -5475.364
$ws.Range("H91").Value = 3711
$ws.Range("J91").Value = 4663.364
$ws.Range("L91").Value = 4663.364
$ws.Range("N91").Value = -7471.364
$ws.Range("H122").Value = 1133.4595
$ws.Range("I122").Value = 943.7857
$ws.Range("J122").Value = 1723.5555
$ws.Range("K122").Value = 2831.3571
$ws.Range("L122").Value = 5170.666499999999
$ws.Range("M122").Value = -381.3571000000002
$ws.Range("N122").Value = -10070.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 54235.668
$ws.Range("I134").Value = 59497.316
$ws.Range("K134").Value = 178491.948
$ws.Range("M134").Value = -175956.948

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2224.2727
$ws.Range("I16").Value = 1716.5
$ws.Range("J16").Value = 3578.3333
$ws.Range("K16").Value = 1716.5
$ws.Range("L16").Value = 3578.3333
$ws.Range("M16").Value = -1429.5
$ws.Range("N16").Value = -4152.3333
$ws.Range("H31").Value = 1216.7106
$ws.Range("I31").Value = 994.4211
$ws.Range("J31").Value = 1439
$ws.Range("K31").Value = 994.4211
$ws.Range("L31").Value = 1439
$ws.Range("M31").Value = -699.4211
$ws.Range("N31").Value = -2029
$ws.Range("H34").Value = 1216.7106
$ws.Range("I34").Value = 994.4211
$ws.Range("J34").Value = 1439
$ws.Range("K34").Value = 994.4211
$ws.Range("L34").Value = 1439
$ws.Range("M34").Value = -792.4211
$ws.Range("N34").Value = -1843
$ws.Range("H58").Value = 3248.611
$ws.Range("J58").Value = 4057.3044
$ws.Range("L58").Value = 4057.3044
$ws.Range("N58").Value = -4463.3044
$ws.Range("H99").Value = 1550
$ws.Range("I99").Value = 1414.2858
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1414.2858
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 83.71419999999989
$ws.Range("N99").Value = -5496
$ws.Range("H113").Value = 2224.2727
$ws.Range("I113").Value = 1716.5
$ws.Range("J113").Value = 3578.3333
$ws.Range("K113").Value = 1716.5
$ws.Range("L113").Value = 3578.3333
$ws.Range("M113").Value = 453.5
$ws.Range("N113").Value = -7918.3333
$ws.Range("H126").Value = 1550
$ws.Range("I126").Value = 1414.2858
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4242.857400000001
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -1772.857400000001
$ws.Range("N126").Value = -12440
$ws.Range("H136").Value = 3248.611
$ws.Range("J136").Value = 4057.3044
$ws.Range("L136").Value = 12171.9132
$ws.Range("N136").Value = -17271.9132

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 680
$ws.Range("J55").Value = 970
$ws.Range("L55").Value = 2910
$ws.Range("N55").Value = -3264
$ws.Range("H68").Value = 927.6495
$ws.Range("I68").Value = 681.0577
$ws.Range("J68").Value = 1212.6
$ws.Range("K68").Value = 2043.1731
$ws.Range("L68").Value = 3637.8
$ws.Range("M68").Value = -1232.1731
$ws.Range("N68").Value = -5259.799999999999
$ws.Range("H71").Value = 927.6495
$ws.Range("I71").Value = 681.0577
$ws.Range("J71").Value = 1212.6
$ws.Range("K71").Value = 6129.5193
$ws.Range("L71").Value = 10913.4
$ws.Range("M71").Value = -2073.5193
$ws.Range("N71").Value = -19025.4
$ws.Range("H105").Value = 908000000
$ws.Range("J105").Value = 908000000
$ws.Range("L105").Value = 2724000000
$ws.Range("N105").Value = -2724005242
$ws.Range("H110").Value = 1027
$ws.Range("I110").Value = 1027
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 3081
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1009
$ws.Range("N110").ClearContents()
$ws.Range("H113").Value = 475.60974
$ws.Range("I113").Value = 468.3913
$ws.Range("J113").Value = 484.83334
$ws.Range("K113").Value = 1405.1739
$ws.Range("L113").Value = 1454.50002
$ws.Range("M113").Value = 764.8261
$ws.Range("N113").Value = -5794.500019999999
$ws.Range("H131").Value = 1614888.1
$ws.Range("I131").Value = 1445.9
$ws.Range("J131").Value = 1925165.4
$ws.Range("K131").Value = 4337.700000000001
$ws.Range("L131").Value = 5775496.199999999
$ws.Range("M131").Value = 702.2999999999993
$ws.Range("N131").Value = -5785576.199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 59229.047
$ws.Range("I80").Value = 2166.6667
$ws.Range("J80").Value = 135312.22
$ws.Range("K80").Value = 2166.6667
$ws.Range("L80").Value = 135312.22
$ws.Range("M80").Value = -1168.6667
$ws.Range("N80").Value = -137308.22
$ws.Range("H83").Value = 59229.047
$ws.Range("I83").Value = 2166.6667
$ws.Range("J83").Value = 135312.22
$ws.Range("K83").Value = 10833.3335
$ws.Range("L83").Value = 676561.1
$ws.Range("M83").Value = -5841.333500000001
$ws.Range("N83").Value = -686545.1
$ws.Range("H126").Value = 1250
$ws.Range("I126").Value = 1000
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 3000
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -530
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 320.25925
$ws.Range("I22").Value = 335.5625
$ws.Range("J22").Value = 298
$ws.Range("K22").Value = 335.5625
$ws.Range("L22").Value = 298
$ws.Range("M22").Value = -40.5625
$ws.Range("N22").Value = -888
$ws.Range("H27").Value = 320.25925
$ws.Range("I27").Value = 335.5625
$ws.Range("J27").Value = 298
$ws.Range("K27").Value = 335.5625
$ws.Range("L27").Value = 298
$ws.Range("M27").Value = -228.5625
$ws.Range("N27").Value = -512
$ws.Range("H40").Value = 2032.6154
$ws.Range("I40").Value = 1552.4
$ws.Range("J40").Value = 3633.3333
$ws.Range("K40").Value = 1552.4
$ws.Range("L40").Value = 3633.3333
$ws.Range("M40").Value = -1416.4
$ws.Range("N40").Value = -3905.3333
$ws.Range("H122").Value = 14863
$ws.Range("I122").Value = 26726
$ws.Range("K122").Value = 80178
$ws.Range("M122").Value = -77728
$ws.Range("H132").Value = 7290.913
$ws.Range("I132").Value = 12776.2
$ws.Range("J132").Value = 3071.4614
$ws.Range("K132").Value = 38328.60000000001
$ws.Range("L132").Value = 9214.3842
$ws.Range("M132").Value = -35798.60000000001
$ws.Range("N132").Value = -14274.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2100.1428
$ws.Range("I126").Value = 2332.6667
$ws.Range("J126").Value = 705
$ws.Range("K126").Value = 6998.000100000001
$ws.Range("L126").Value = 2115
$ws.Range("M126").Value = -4528.000100000001
$ws.Range("N126").Value = -7055
$ws.Range("H132").Value = 5969.674
$ws.Range("I132").Value = 6673.4116
$ws.Range("J132").Value = 3975.75
$ws.Range("K132").Value = 20020.2348
$ws.Range("L132").Value = 11927.25
$ws.Range("M132").Value = -17490.2348
$ws.Range("N132").Value = -16987.25
